$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 6 (the last data row) - also shrinks the used dimension to A1:AH5
$ws.Rows.Item(6).Delete()

# 2) Apply "custom accuracy" rounding (2 decimal places) to the remaining last
#    data row (row 5), replacing the raw 3-decimal sensor readings with their
#    rounded values.
$ws.Range("B5").Value = 14.75
$ws.Range("C5").Value = 11.02
$ws.Range("D5").Value = 1.02
$ws.Range("E5").Value = 32.39
$ws.Range("F5").Value = 26.24
$ws.Range("G5").Value = 11.34
$ws.Range("H5").Value = 44.88
$ws.Range("I5").Value = 17.99
$ws.Range("J5").Value = 8.08
$ws.Range("K5").Value = 11.57
$ws.Range("L5").Value = 13.13
$ws.Range("M5").Value = 13.73
$ws.Range("N5").Value = 3.93
$ws.Range("O5").Value = 11.65
$ws.Range("P5").Value = 16.47
$ws.Range("Q5").Value = 9.960000000000001
$ws.Range("R5").Value = 0.65
$ws.Range("S5").Value = 0.49
$ws.Range("T5").Value = 170.07
$ws.Range("U5").Value = 32.58
$ws.Range("V5").Value = 10.76
$ws.Range("W5").Value = 21.76
$ws.Range("X5").Value = 11.5
$ws.Range("Y5").Value = 1.56
$ws.Range("Z5").Value = 22.37
$ws.Range("AA5").Value = 9.5
$ws.Range("AB5").Value = 8.48
$ws.Range("AC5").Value = 9.949999999999999
$ws.Range("AD5").Value = 13.66
$ws.Range("AE5").Value = 0.5
$ws.Range("AF5").Value = 40.92
$ws.Range("AG5").Value = 6
$ws.Range("AH5").Value = 13.45

# 3) Narrow a handful of columns (C, G, K, L, V, X -> indices 3,7,11,12,22,24)
#    by one character (OOXML width 8 -> 7; Excel's ColumnWidth COM property
#    is offset by -0.83 from the stored OOXML column width on this sheet's
#    font).
$narrowCols = @(3, 7, 11, 12, 22, 24)
foreach ($col in $narrowCols) {
    $ws.Columns.Item($col).ColumnWidth = 6.17
}
